$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 79, shifting existing rows 79-155 down to 80-156.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new data record.
$ws.Cells.Item(79, 1).Value = 5
$ws.Cells.Item(79, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(79, 3).Value = "Maule"
$ws.Cells.Item(79, 4).Value = 44880
$ws.Cells.Item(79, 5).Value = 7
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100108
$ws.Cells.Item(79, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(79, 9).Value = 100108002
$ws.Cells.Item(79, 10).Value = "Mango"
$ws.Cells.Item(79, 11).Value = "Sin especificar"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 248
$ws.Cells.Item(79, 14).Value = 8000
$ws.Cells.Item(79, 15).Value = 8000
$ws.Cells.Item(79, 16).Value = 8000
$ws.Cells.Item(79, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(79, 18).Value = "Brasil"
$ws.Cells.Item(79, 19).Value = 2000
$ws.Cells.Item(79, 20).Value = 4
